# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> bound to the (one and only) slide master,
#                             currently the "Integral" / "Red Violet" theme
#   ppt/theme/theme2.xml  -> bound only to the notes master,
#                             currently the default "Office Theme"
#
# The authored change swaps the two themes' contents: the slide master
# (and therefore every slide) should end up using the plain "Office"
# colour scheme, while the notes master ends up with the old
# "Integral" / "Red Violet" palette.
#
# PowerPoint's object model doesn't give us a raw "swap part bytes"
# verb, so we reproduce the same end effect the way an author driving
# the Design/Variants UI would: by pushing the target RGB values onto
# the live theme colour scheme, slot by slot, through
# ThemeColorScheme.Colors(i).RGB.

$p = $ppt.ActivePresentation

function Pack-RGB([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette = the stock Office theme colours (what theme1.xml
# should contain after the edit). Order matches the standard
# ThemeColorScheme.Colors index layout:
#  1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#  8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeColors = @(
    (Pack-RGB 0x00 0x00 0x00), # dk1      000000
    (Pack-RGB 0xFF 0xFF 0xFF), # lt1      FFFFFF
    (Pack-RGB 0x44 0x54 0x6A), # dk2      44546A
    (Pack-RGB 0xE7 0xE6 0xE6), # lt2      E7E6E6
    (Pack-RGB 0x5B 0x9B 0xD5), # accent1  5B9BD5
    (Pack-RGB 0xED 0x7D 0x31), # accent2  ED7D31
    (Pack-RGB 0xA5 0xA5 0xA5), # accent3  A5A5A5
    (Pack-RGB 0xFF 0xC0 0x00), # accent4  FFC000
    (Pack-RGB 0x44 0x72 0xC4), # accent5  4472C4
    (Pack-RGB 0x70 0xAD 0x47), # accent6  70AD47
    (Pack-RGB 0x05 0x63 0xC1), # hlink    0563C1
    (Pack-RGB 0x95 0x4F 0x72)  # folHlink 954F72
)

$tcs = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $tcs.Colors($i).RGB = $officeColors[$i - 1]
}
